$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 3-25 with new data (row 1 header, row 2 unchanged)
$ws.Cells.Item(3, 1).Value = 591
$ws.Cells.Item(3, 2).Value = 'Doxel.ai'
$ws.Cells.Item(3, 3).Value = 'CS1 Doxel - Enterprise AE Northeast U.S'
$ws.Cells.Item(3, 4).Value = 'Matt Crandley'
$ws.Cells.Item(3, 5).Value = '2nd Interview'

$ws.Cells.Item(4, 1).Value = 694
$ws.Cells.Item(4, 2).Value = 'Chronosphere'
$ws.Cells.Item(4, 3).Value = 'CS1 Chronosphere - Enterprise AE'
$ws.Cells.Item(4, 4).Value = 'Frank Casper'
$ws.Cells.Item(4, 5).Value = '1st Interview'

$ws.Cells.Item(5, 1).Value = 696
$ws.Cells.Item(5, 2).Value = 'Cognition AI'
$ws.Cells.Item(5, 3).Value = 'Founding EMEA AE'
$ws.Cells.Item(5, 4).Value = 'Felix Völker'
$ws.Cells.Item(5, 5).Value = 'CV Sent'

$ws.Cells.Item(6, 1).Value = 715
$ws.Cells.Item(6, 2).Value = 'Honeycomb'
$ws.Cells.Item(6, 3).Value = 'Strat AE U.S x4'
$ws.Cells.Item(6, 4).Value = 'Brian Saverino'
$ws.Cells.Item(6, 5).Value = '4th Interview'

$ws.Cells.Item(7, 1).Value = 731
$ws.Cells.Item(7, 2).Value = 'Oscilar'
$ws.Cells.Item(7, 3).Value = 'Enterprise AE x5'
$ws.Cells.Item(7, 4).Value = 'Greg Muender'
$ws.Cells.Item(7, 5).Value = '4th Interview'

$ws.Cells.Item(8, 1).Value = 740
$ws.Cells.Item(8, 2).Value = 'Axion Ray'
$ws.Cells.Item(8, 3).Value = 'Enterprise Account Executive (East)'
$ws.Cells.Item(8, 4).Value = 'Lindsay St.Cin'
$ws.Cells.Item(8, 5).Value = '1st Interview'

$ws.Cells.Item(9, 1).Value = 743
$ws.Cells.Item(9, 2).Value = 'Conduct AI'
$ws.Cells.Item(9, 3).Value = 'CS1 Conduct AI - Enterprise Account Executive UK (German speaking)'
$ws.Cells.Item(9, 4).Value = 'Leonard Friederich'
$ws.Cells.Item(9, 5).Value = '3rd Interview'

$ws.Cells.Item(10, 1).Value = 746
$ws.Cells.Item(10, 2).Value = 'LaunchDarkly'
$ws.Cells.Item(10, 3).Value = 'Enterprise Account Executive (Germany)'
$ws.Cells.Item(10, 4).Value = 'Aron Kraft'
$ws.Cells.Item(10, 5).Value = '1st Interview'

$ws.Cells.Item(11, 1).Value = 746
$ws.Cells.Item(11, 2).Value = 'LaunchDarkly'
$ws.Cells.Item(11, 3).Value = 'Enterprise Account Executive (Germany)'
$ws.Cells.Item(11, 4).Value = 'Florian Werner'
$ws.Cells.Item(11, 5).Value = '1st Interview'

$ws.Cells.Item(12, 1).Value = 746
$ws.Cells.Item(12, 2).Value = 'LaunchDarkly'
$ws.Cells.Item(12, 3).Value = 'Enterprise Account Executive (Germany)'
$ws.Cells.Item(12, 4).Value = 'Maximilian May'
$ws.Cells.Item(12, 5).Value = '1st Interview'

$ws.Cells.Item(13, 1).Value = 801
$ws.Cells.Item(13, 2).Value = 'Redwood Software'
$ws.Cells.Item(13, 3).Value = 'Redwood AE Germany x 3'
$ws.Cells.Item(13, 4).Value = 'Thomas Schaeffer'
$ws.Cells.Item(13, 5).Value = 'CV Sent'

$ws.Cells.Item(14, 1).Value = 801
$ws.Cells.Item(14, 2).Value = 'Redwood Software'
$ws.Cells.Item(14, 3).Value = 'Redwood AE Germany x 3'
$ws.Cells.Item(14, 4).Value = 'Leonard Friederich'
$ws.Cells.Item(14, 5).Value = '4th Interview'

$ws.Cells.Item(15, 1).Value = 801
$ws.Cells.Item(15, 2).Value = 'Redwood Software'
$ws.Cells.Item(15, 3).Value = 'Redwood AE Germany x 3'
$ws.Cells.Item(15, 4).Value = 'Felix Völker'
$ws.Cells.Item(15, 5).Value = 'CV Sent'

$ws.Cells.Item(16, 1).Value = 801
$ws.Cells.Item(16, 2).Value = 'Redwood Software'
$ws.Cells.Item(16, 3).Value = 'Redwood AE Germany x 3'
$ws.Cells.Item(16, 4).Value = 'Konstantin Melzer'
$ws.Cells.Item(16, 5).Value = 'CV Sent'

$ws.Cells.Item(17, 1).Value = 810
$ws.Cells.Item(17, 2).Value = 'groundcover'
$ws.Cells.Item(17, 3).Value = 'Mid-Market AE (Observability in NYC, Boston, Denver, SF)'
$ws.Cells.Item(17, 4).Value = 'Alex Andrei'
$ws.Cells.Item(17, 5).Value = '1st Interview'

$ws.Cells.Item(18, 1).Value = 833
$ws.Cells.Item(18, 2).Value = 'Blockaid'
$ws.Cells.Item(18, 3).Value = 'SDR Manager'
$ws.Cells.Item(18, 4).Value = 'Jacques Lane'
$ws.Cells.Item(18, 5).Value = 'CV Sent'

$ws.Cells.Item(19, 1).Value = 839
$ws.Cells.Item(19, 2).Value = 'Pigment'
$ws.Cells.Item(19, 3).Value = 'Enterprise AE (California)'
$ws.Cells.Item(19, 4).Value = 'Alexander Wise'
$ws.Cells.Item(19, 5).Value = 'CV Sent'

$ws.Cells.Item(20, 1).Value = 839
$ws.Cells.Item(20, 2).Value = 'Pigment'
$ws.Cells.Item(20, 3).Value = 'Enterprise AE (California)'
$ws.Cells.Item(20, 4).Value = 'Colin Glen'
$ws.Cells.Item(20, 5).Value = 'CV Sent'

$ws.Cells.Item(21, 1).Value = 839
$ws.Cells.Item(21, 2).Value = 'Pigment'
$ws.Cells.Item(21, 3).Value = 'Enterprise AE (California)'
$ws.Cells.Item(21, 4).Value = 'Cris Castillo'
$ws.Cells.Item(21, 5).Value = '2nd Interview'

$ws.Cells.Item(22, 1).Value = 847
$ws.Cells.Item(22, 2).Value = 'Simile.ai'
$ws.Cells.Item(22, 3).Value = 'CS1 Simile.ai - Enterprise AE x2'
$ws.Cells.Item(22, 4).Value = 'ROMIT MIRCHANDANI'
$ws.Cells.Item(22, 5).Value = '2nd Interview'

$ws.Cells.Item(23, 1).Value = 847
$ws.Cells.Item(23, 2).Value = 'Simile.ai'
$ws.Cells.Item(23, 3).Value = 'CS1 Simile.ai - Enterprise AE x2'
$ws.Cells.Item(23, 4).Value = 'Ali Hussain'
$ws.Cells.Item(23, 5).Value = 'CV Sent'

$ws.Cells.Item(24, 1).Value = 847
$ws.Cells.Item(24, 2).Value = 'Simile.ai'
$ws.Cells.Item(24, 3).Value = 'CS1 Simile.ai - Enterprise AE x2'
$ws.Cells.Item(24, 4).Value = 'Danny Wymer'
$ws.Cells.Item(24, 5).Value = '3rd Interview'

$ws.Cells.Item(25, 1).Value = 850
$ws.Cells.Item(25, 2).Value = 'Laurel'
$ws.Cells.Item(25, 3).Value = 'Enterprise AE New York'
$ws.Cells.Item(25, 4).Value = 'Matt Crandley'
$ws.Cells.Item(25, 5).Value = '1st Interview'

# Clear now-unused rows 26-35 (shrinks used range to A1:E25)
$ws.Range("A26:E35").ClearContents()
